$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.383.64"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "2.794.85"
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "347.25"
$ws.Range("E5").Value = "  +4.01%  "

$ws.Range("D6").Value = "115.96"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  +3.49%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  +2.33%  "

$ws.Range("D10").Value = "42.53"
$ws.Range("E10").Value = "  +2.16%  "

$ws.Range("D11").Value = "0.0861"
$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("D12").Value = "20.03"
$ws.Range("E12").Value = "  -1.20%  "

$ws.Range("E13").Value = "  +1.76%  "

$ws.Range("D14").Value = "7.90"
$ws.Range("E14").Value = "  +3.40%  "

$ws.Range("D15").Value = "3.228.58"
$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("D16").Value = "2.777.72"
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("D17").Value = "0.893"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "52.199.16"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("E19").Value = "  +5.93%  "

$ws.Range("D20").Value = "7.29"
$ws.Range("E20").Value = "  +5.91%  "

$ws.Range("D21").Value = "13.38"
$ws.Range("E21").Value = "  -4.03%  "

$ws.Range("D22").Value = "0.0₃0982"
$ws.Range("E22").Value = "  +1.93%  "

$ws.Range("D23").Value = "270.19"
$ws.Range("E23").Value = "  -3.20%  "

$ws.Range("D24").Value = "70.06"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("E25").Value = "  +2.84%  "

$ws.Range("D26").Value = "26.85"
$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  -1.45%  "

$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "34.73"
$ws.Range("E31").Value = "  -3.34%  "

$ws.Range("B32").Value = "VeChain"
$ws.Range("C32").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D32").Value = "0.0461"
$ws.Range("E32").Value = "  +31.70%  "

$ws.Range("D33").Value = "50.14"
$ws.Range("E33").Value = "  -0.58%  "

$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").Value = "0.0837"
$ws.Range("E35").Value = "  +0.95%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "4.95"
$ws.Range("E38").Value = "  -1.14%  "

$ws.Range("D39").Value = "18.68"
$ws.Range("E39").Value = "  -4.26%  "

$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  +10.67%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "23.54"
$ws.Range("E42").Value = "  -1.27%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "126.92"
$ws.Range("E43").Value = "  -2.08%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.115"
$ws.Range("E44").Value = "  +1.88%  "

$ws.Range("D45").Value = "2.29"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("E46").Value = "  -2.42%  "

$ws.Range("D47").Value = "2.059.85"
$ws.Range("E47").Value = "  -2.30%  "

$ws.Range("E48").Value = "  +3.35%  "

$ws.Range("D49").Value = "0.967"
$ws.Range("E49").Value = "  +11.11%  "

$ws.Range("D50").Value = "5.61"
$ws.Range("E50").Value = "  -0.27%  "

$ws.Range("D51").Value = "8.99"
$ws.Range("E51").Value = "  -0.82%  "
